$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 76; this shifts existing rows 76..189 down to 77..190
$ws.Rows.Item(76).Insert()

# Populate the newly inserted row 76 with the new data record
$ws.Cells.Item(76, 1).Value = 11
$ws.Cells.Item(76, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(76, 3).Value = "Bíobío"
$ws.Cells.Item(76, 4).Value = 44771
$ws.Cells.Item(76, 5).Value = 8
$ws.Cells.Item(76, 6).Value = 100112003
$ws.Cells.Item(76, 7).Value = "Ajo"
$ws.Cells.Item(76, 8).Value = "Chino"
$ws.Cells.Item(76, 9).Value = "Primera"
$ws.Cells.Item(76, 10).Value = 400
$ws.Cells.Item(76, 11).Value = 24000
$ws.Cells.Item(76, 12).Value = 25000
$ws.Cells.Item(76, 13).Value = 24500
$ws.Cells.Item(76, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(76, 15).Value = "China"
$ws.Cells.Item(76, 16).Value = 2450
$ws.Cells.Item(76, 17).Value = 10
$ws.Cells.Item(76, 18).Value = "Hortaliza"
